# Auto - Update data with bot!
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 26: update title
$ws.Range("D26").Value = "ai plus(est soft)"

# Row 44: update title and link
$ws.Range("D44").Value = "SK 스퀘어 포트폴리오 분석 (1) - 코빗"
$ws.Range("E44").Value = "https://engineering-ladder.tistory.com/112"

# Row 50: update title and link
$ws.Range("D50").Value = "가장 쉬운 PyMC3 예제"
$ws.Range("E50").Value = "http://incredible.egloos.com/7534316"
